$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The New India Assurance row (old row 5) is folded away; the data that used to
# live in rows 2-4 (company "3", Cholamandalam, ICICI Lombard) is replaced with a
# refreshed 3-company table (index "2", ICICI Lombard, The New India Assurance).
# Deleting row 5 shifts nothing below it and Excel recomputes the sheet dimension
# to A1:AQ4 automatically.
$ws.Rows("5:5").Delete()

# --- Row 2 ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("B2").ClearFormats()
$ws.Range("D2").Value = [double]"0.16"
$ws.Range("E2").Value = [double]"0.198"
$ws.Range("F2").Value = [double]"0.1185"
$ws.Range("G2").Value = [double]"0.07990938645638793"
$ws.Range("H2").Value = [double]"0.07990938645638793"
$ws.Range("I2").Value = [double]"0.08727661855357467"
$ws.Range("J2").Value = [double]"0.06870250687117567"
$ws.Range("K2").Value = [double]"391.4"
$ws.Range("L2").Value = [double]"0.0676834751331535"
$ws.Range("M2").Value = [double]"0.244"
$ws.Range("N2").Value = [double]"1.972338757264916e-05"
$ws.Range("O2").Value = [double]"0.000623403168114461"
$ws.Range("P2").Value = [double]"0.244"
$ws.Range("Q2").Value = [double]"1.972338757264916e-05"
$ws.Range("R2").Value = [double]"0.000623403168114461"
$ws.Range("U2").Value = [double]"1587.15"
$ws.Range("V2").Value = [double]"0.128294977811189"
$ws.Range("W2").Value = [double]"0.1337420777421839"
$ws.Range("X2").Value = [double]"0.05426910984075353"
$ws.Range("Y2").Value = [double]"0.07947296790143038"
$ws.Range("Z2").Value = [double]"1.220848966229092"
$ws.Range("AA2").Value = [double]"0.1377060874950877"
$ws.Range("AB2").Value = [double]"0.05417451168483343"
$ws.Range("AC2").Value = [double]"0.08353157581025424"
$ws.Range("AD2").Value = [double]"66"
$ws.Range("AE2").Value = [double]"0.003851141942463297"
$ws.Range("AF2").Value = [double]"66.00385114194246"
$ws.Range("AG2").Value = [double]"-1521.146148858058"
$ws.Range("AH2").Value = [double]"0.00530701133736068"
$ws.Range("AI2").Value = [double]"0.01200920741849313"
$ws.Range("AJ2").Value = [double]"-0.1401983980510626"
$ws.Range("AK2").Value = [double]"-0.3891440540833393"
$ws.Range("AL2").Value = [double]"5.538"
$ws.Range("AM2").Value = [double]"5.538"
$ws.Range("AN2").Value = [double]"0.1282788860728002"
$ws.Range("AO2").Value = [double]"91.13398338750451"
$ws.Range("AP2").Value = [double]"-2.95652929590063"
$ws.Range("AQ2").Value = [double]"91.13398338750451"

# --- Row 3 ---
$ws.Range("B3").Value = "ICICI Lombard General Insurance Company Limited (BSE:540716)"
$ws.Range("D3").Value = [double]"0.16"
$ws.Range("E3").Value = [double]"0.198"
$ws.Range("F3").Value = [double]"0.187"
$ws.Range("G3").Value = [double]"0.1468103883005547"
$ws.Range("H3").Value = [double]"0.1468103883005547"
$ws.Range("I3").Value = [double]"0.1601760147324833"
$ws.Range("J3").Value = [double]"0.1205625917341272"
$ws.Range("K3").Value = [double]"189"
$ws.Range("L3").Value = [double]"0.1191376701966717"
$ws.Range("M3").Value = [double]"-0"
$ws.Range("N3").Value = [double]"-0"
$ws.Range("O3").Value = [double]"-0"
$ws.Range("P3").Value = [double]"-0"
$ws.Range("Q3").Value = [double]"-0"
$ws.Range("R3").Value = [double]"-0"
$ws.Range("U3").Value = [double]"3.65"
$ws.Range("V3").Value = [double]"0.0003856109027520997"
$ws.Range("W3").Value = [double]"0.2282608695652174"
$ws.Range("X3").Value = [double]"0.05437859352681861"
$ws.Range("Y3").Value = [double]"0.1738822760383988"
$ws.Range("Z3").Value = [double]"1.843785331614338"
$ws.Range("AA3").Value = [double]"0.2222915381807918"
$ws.Range("AB3").Value = [double]"0.05418939721497843"
$ws.Range("AC3").Value = [double]"0.1681021409658134"
$ws.Range("AD3").Value = [double]"66"
$ws.Range("AE3").Value = [double]"0.003851141942463297"
$ws.Range("AF3").Value = [double]"66.00385114194246"
$ws.Range("AG3").Value = [double]"62.35385114194246"
$ws.Range("AH3").Value = [double]"0.006924809785817243"
$ws.Range("AI3").Value = [double]"0.06339795119338541"
$ws.Range("AJ3").Value = [double]"0.006544375272346265"
$ws.Range("AK3").Value = [double]"0.06010277090716715"
$ws.Range("AL3").Value = [double]"5.45"
$ws.Range("AM3").Value = [double]"5.45"
$ws.Range("AN3").Value = [double]"0.2568053415511042"
$ws.Range("AO3").Value = [double]"46.62385321100917"
$ws.Range("AP3").Value = [double]"0.2426182127202007"
$ws.Range("AQ3").Value = [double]"46.62385321100917"
$ws.Range("T3").ClearContents()

# --- Row 4 ---
$ws.Range("B4").Value = "The New India Assurance Company Limited (BSE:540769)"
$ws.Range("F4").Value = [double]"0.05"
$ws.Range("G4").Value = [double]"0.05461824420932228"
$ws.Range("H4").Value = [double]"0.05461824420932228"
$ws.Range("I4").Value = [double]"0.0597178533981508"
$ws.Range("J4").Value = [double]"0.0490686122543641"
$ws.Range("K4").Value = [double]"202.4"
$ws.Range("L4").Value = [double]"0.04823181774854638"
$ws.Range("M4").Value = [double]"0.244"
$ws.Range("N4").Value = [double]"8.397577092511013e-05"
$ws.Range("O4").Value = [double]"0.001205533596837945"
$ws.Range("P4").Value = [double]"0.244"
$ws.Range("Q4").Value = [double]"8.397577092511013e-05"
$ws.Range("R4").Value = [double]"0.001205533596837945"
$ws.Range("U4").Value = [double]"1583.5"
$ws.Range("V4").Value = [double]"0.544982103524229"
$ws.Range("W4").Value = [double]"0.03922328591915043"
$ws.Range("X4").Value = [double]"0.05415962615468844"
$ws.Range("Y4").Value = [double]"-0.01493634023553802"
$ws.Range("Z4").Value = [double]"1.08257874777494"
$ws.Range("AA4").Value = [double]"0.05312063680938357"
$ws.Range("AB4").Value = [double]"0.05415962615468844"
$ws.Range("AC4").Value = [double]"-0.001038989345304879"
$ws.Range("AD4").Value = [double]"0"
$ws.Range("AE4").Value = [double]"0"
$ws.Range("AF4").Value = [double]"0"
$ws.Range("AG4").Value = [double]"-1583.5"
$ws.Range("AH4").Value = [double]"0"
$ws.Range("AI4").Value = [double]"0"
$ws.Range("AJ4").Value = [double]"-1.197715755237879"
$ws.Range("AK4").Value = [double]"-0.551453943931743"
$ws.Range("AL4").Value = [double]"0.08799999999999999"
$ws.Range("AM4").Value = [double]"0.08799999999999999"
$ws.Range("AN4").Value = [double]"0"
$ws.Range("AO4").Value = [double]"2847.727272727273"
$ws.Range("AP4").Value = [double]"-6.149514563106796"
$ws.Range("AQ4").Value = [double]"2847.727272727273"
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
